$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-prioritized backlog items following sprint 1 demo.
# Column A = Priority (re-estimated), Column C/D = backlog item text (re-ordered).

$ws.Range("A2").Value = 0.5
$ws.Range("A3").Value = 0.5
$ws.Range("A4").Value = 0.5
$ws.Range("A8").Value = 1.5
$ws.Range("A9").Value = 1.5

$ws.Range("C4").Value = "View an overview of my trip"
$ws.Range("D4").Value = "I can see my trip plans at a glance"

$ws.Range("C5").Value = "Remove Waypoints from my trip"
$ws.Range("D5").Value = "I can clear up space in my trip and/or remove places I no longer want to go to"

$ws.Range("C6").Value = "Add Transportation to my trip"
$ws.Range("D6").Value = "I can record how I plan to travel to/from Waypoints"

$ws.Range("C7").Value = "Remove Transportation from my trip"
$ws.Range("D7").Value = "I can remove details on how I plan to travel to/from Waypoints if it is not needed (removed Waypoint, travel is trivial, or no longer want to use specified travel option and don't have a replacement travel option)"

$ws.Range("C8").Value = "Add Lodging to a trip"
$ws.Range("D8").Value = "I can specify where I will be staying during a specified time period of my trip"

$ws.Range("C9").Value = "Remove Lodging from a trip"
$ws.Range("D9").Value = "I can remove Lodging that I will no longer be using, but have not found or do not want to specify replacement lodging"

# View state updates: zoom in to 160% and leave the A8:C9 block selected
# (as it was after reviewing the re-prioritized rows during the demo).
$excel.ActiveWindow.Zoom = 160
$ws.Range("A8:C9").Select()
